$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (id, speaker_variant) pairs for rows 2-9, with is_prefered (D) cleared
$rows = @(
    @{ Row = 2; B = "#kalofisi";  C = "Kalofisi" },
    @{ Row = 3; B = "#jonker";    C = "Jonker" },
    @{ Row = 4; B = "#schout";    C = "Schout" },
    @{ Row = 5; B = "#karlofisi"; C = "Karlofisi" },
    @{ Row = 6; B = "#augestyn";  C = "Augestyn" },
    @{ Row = 7; B = "#meindert";  C = "Meindert" },
    @{ Row = 8; B = "#barent";    C = "Barent" },
    @{ Row = 9; B = "#izabel";    C = "Izabel" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $null
}
